$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 15250
$ws.Range("I16").Value = 29500
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 29500
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -29270
$ws.Range("N16").Value = -1460
$ws.Range("H17").Value = 1490.7778
$ws.Range("J17").Value = 1490.7778
$ws.Range("L17").Value = 4472.3334
$ws.Range("N17").Value = -4808.3334
$ws.Range("H21").Value = 15254
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 15254
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 15254
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -16190
$ws.Range("H23").Value = 15254
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 15254
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 15254
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -15722
$ws.Range("H32").Value = 1125
$ws.Range("J32").Value = 1125
$ws.Range("L32").Value = 1125
$ws.Range("N32").Value = -1777
$ws.Range("H135").Value = 894.6667
$ws.Range("I135").Value = 842
$ws.Range("K135").Value = 7578
$ws.Range("M135").Value = -5043
$ws.Range("H138").Value = 12536.132
$ws.Range("I138").Value = 5000
$ws.Range("J138").Value = 13182.086
$ws.Range("K138").Value = 15000
$ws.Range("L138").Value = 39546.258
$ws.Range("M138").Value = -9860
$ws.Range("N138").Value = -49826.258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1370
$ws.Range("I2").Value = 1370.0714
$ws.Range("K2").Value = 1370.0714
$ws.Range("M2").Value = -1257.0714
$ws.Range("H110").Value = 1047.25
$ws.Range("I110").Value = 1058.5555
$ws.Range("J110").Value = 1013.3333
$ws.Range("K110").Value = 1058.5555
$ws.Range("L110").Value = 1013.3333
$ws.Range("M110").Value = 986.4445000000001
$ws.Range("N110").Value = -5103.3333
$ws.Range("H116").Value = 1370
$ws.Range("I116").Value = 1370.0714
$ws.Range("K116").Value = 1370.0714
$ws.Range("M116").Value = 923.9286
$ws.Range("H122").Value = 1523.1818
$ws.Range("I122").Value = 1425.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4276.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1826.5
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1370
$ws.Range("I3").Value = 1370.0714
$ws.Range("K3").Value = 1370.0714
$ws.Range("M3").Value = -1256.0714
$ws.Range("H86").Value = 6916.5
$ws.Range("I86").Value = 3633.1667
$ws.Range("K86").Value = 3633.1667
$ws.Range("M86").Value = -2510.1667
$ws.Range("H89").Value = 6916.5
$ws.Range("I89").Value = 3633.1667
$ws.Range("K89").Value = 18165.8335
$ws.Range("M89").Value = -12549.8335
$ws.Range("H94").Value = 5314
$ws.Range("I94").Value = 3212.75
$ws.Range("J94").Value = 10357
$ws.Range("K94").Value = 3212.75
$ws.Range("L94").Value = 10357
$ws.Range("M94").Value = -2761.75
$ws.Range("N94").Value = -11259
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H140").Value = 51400
$ws.Range("J140").Value = 51400
$ws.Range("L140").Value = 51400
$ws.Range("N140").Value = -61760

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1490.7858
$ws.Range("I7").Value = 2252.3333
$ws.Range("K7").Value = 2252.3333
$ws.Range("M7").Value = -2139.3333
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -407
$ws.Range("H58").Value = 4298.1665
$ws.Range("I58").Value = 4157.8
$ws.Range("K58").Value = 4157.8
$ws.Range("M58").Value = -3954.8
$ws.Range("H94").Value = 4170.1665
$ws.Range("I94").Value = 3001.4
$ws.Range("K94").Value = 3001.4
$ws.Range("M94").Value = -2550.4
$ws.Range("H136").Value = 4298.1665
$ws.Range("I136").Value = 4157.8
$ws.Range("K136").Value = 12473.4
$ws.Range("M136").Value = -9923.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3082.2
$ws.Range("I102").Value = 3082.2
$ws.Range("K102").Value = 3082.2
$ws.Range("M102").Value = -1460.2
$ws.Range("H122").Value = 7732.8335
$ws.Range("J122").Value = 10196.667
$ws.Range("L122").Value = 30590.001
$ws.Range("N122").Value = -35490.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1647.75
$ws.Range("I16").Value = 1647.75
$ws.Range("K16").Value = 1647.75
$ws.Range("M16").Value = -1477.75
$ws.Range("H61").Value = 7286935
$ws.Range("I61").Value = 5667866
$ws.Range("K61").Value = 5667866
$ws.Range("M61").Value = -5667664
$ws.Range("H68").Value = 1444
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 1444
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 7286935
$ws.Range("I113").Value = 5667866
$ws.Range("K113").Value = 5667866
$ws.Range("M113").Value = -5665696
$ws.Range("H132").Value = 3639.6667
$ws.Range("I132").Value = 2749.6667
$ws.Range("J132").Value = 4974.6665
$ws.Range("K132").Value = 8249.000100000001
$ws.Range("L132").Value = 14923.9995
$ws.Range("M132").Value = -5719.000100000001
$ws.Range("N132").Value = -19983.9995
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 395734.84
$ws.Range("I2").Value = 140481.8
$ws.Range("J2").Value = 1672000
$ws.Range("K2").Value = 140481.8
$ws.Range("L2").Value = 1672000
$ws.Range("M2").Value = -140369.8
$ws.Range("N2").Value = -1672224
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H74").Value = 34318.2
$ws.Range("J74").Value = 31648
$ws.Range("L74").Value = 31648
$ws.Range("N74").Value = -33520
$ws.Range("H77").Value = 34318.2
$ws.Range("J77").Value = 31648
$ws.Range("L77").Value = 94944
$ws.Range("N77").Value = -104304
$ws.Range("H125").Value = 99998
$ws.Range("J125").Value = 99998
$ws.Range("L125").Value = 99998
$ws.Range("N125").Value = -109838
